$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '274.88'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '22.96'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '6.351'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.06265'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '3.664'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '6.679'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.8314'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.01376'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.1629'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08369'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.03462'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.03106'
$ws.Range('B15').Value = 'ProBitToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.1241'
$ws.Range('E15').Value = '14ProBitTokenPROB'
$ws.Range('B16').Value = 'BitMartToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.09307'
$ws.Range('E16').Value = '15BitMartTokenBMX'
$ws.Range('B17').Value = 'MCDex'
$ws.Range('C17').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.906'
$ws.Range('E17').Value = '16MCDexMCB'
$ws.Range('B18').Value = 'BitForexToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.001642'
$ws.Range('E18').Value = '17BitForexTokenBF'
$ws.Range('B19').Value = 'CoinExToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.04771'
$ws.Range('E19').Value = '18CoinExTokenCET'
$ws.Range('B20').Value = 'TigerCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.006378'
$ws.Range('E20').Value = '19TigerCashTCH'
$ws.Range('B21').Value = 'HotbitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.005687'
$ws.Range('E21').Value = '20HotbitTokenHTB'
$ws.Range('B22').Value = 'BitKan'
$ws.Range('C22').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.001095'
$ws.Range('E22').Value = '21BitKanKAN'
$ws.Range('B23').Value = 'NitroEx'
$ws.Range('C23').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.0001500'
$ws.Range('E23').Value = '22NitroExNTX'
$ws.Range('B24').Value = 'LEO'
$ws.Range('C24').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '3.715'
$ws.Range('E24').Value = '23LEOLEO'
$ws.Range('B25').Value = 'BTSEToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.370'
$ws.Range('E25').Value = '24BTSETokenBTSE'
$ws.Range('B26').Value = 'BitpandaEcosystemToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.3342'
$ws.Range('E26').Value = '25BitpandaEcosystemTokenBEST'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0002680'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.04707'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.007022'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.003349'
$ws.Range('E43').Value = '42CEJICEJIWorstin24h'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00006254'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.8999'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.03052'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.00002200'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.01240'
